# Update "想去人数" (column F) counts across the sheets of the
# 北京-漫展信息 workbook, as generated at commit 456a3b4.
#
# Sheet order in the workbook:
#   1 = 展览      (Exhibitions)
#   2 = 演出      (Performances)
#   3 = 本地生活  (Local life)
#   4 = 全部类型  (All types, aggregates the other three sheets)

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("F5").Value = 318
$ws.Range("F6").Value = 398
$ws.Range("F8").Value = 46
$ws.Range("F9").Value = 513
$ws.Range("F13").Value = 105
$ws.Range("F14").Value = 238
$ws.Range("F15").Value = 34
$ws.Range("F16").Value = 412
$ws.Range("F17").Value = 6620
$ws.Range("F18").Value = 63
$ws.Range("F19").Value = 71
$ws.Range("F21").Value = 7567
$ws.Range("F23").Value = 36
$ws.Range("F24").Value = 3392
$ws.Range("F25").Value = 27
$ws.Range("F26").Value = 1791
$ws.Range("F27").Value = 891
$ws.Range("F28").Value = 4514
$ws.Range("F29").Value = 112
$ws.Range("F32").Value = 217
$ws.Range("F34").Value = 1663
$ws.Range("F36").Value = 163
$ws.Range("F37").Value = 55
$ws.Range("F39").Value = 1201
$ws.Range("F40").Value = 1754

# --- Sheet 2: 演出 ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("F2").Value = 5
$ws.Range("F3").Value = 67
$ws.Range("F4").Value = 50
$ws.Range("F5").Value = 5

# --- Sheet 3: 本地生活 ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("F3").Value = 1223

# --- Sheet 4: 全部类型 ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("F4").Value = 1223
$ws.Range("F7").Value = 318
$ws.Range("F8").Value = 398
$ws.Range("F10").Value = 46
$ws.Range("F11").Value = 513
$ws.Range("F12").Value = 5
$ws.Range("F15").Value = 67
$ws.Range("F16").Value = 105
$ws.Range("F17").Value = 238
$ws.Range("F18").Value = 34
$ws.Range("F19").Value = 412
$ws.Range("F20").Value = 6620
$ws.Range("F21").Value = 63
$ws.Range("F22").Value = 71
$ws.Range("F24").Value = 7567
$ws.Range("F26").Value = 36
$ws.Range("F27").Value = 3392
$ws.Range("F28").Value = 27
$ws.Range("F29").Value = 1791
$ws.Range("F30").Value = 891
$ws.Range("F31").Value = 4514
$ws.Range("F32").Value = 112
$ws.Range("F35").Value = 50
$ws.Range("F36").Value = 217
$ws.Range("F38").Value = 1663
$ws.Range("F40").Value = 163
$ws.Range("F41").Value = 55
$ws.Range("F43").Value = 5
$ws.Range("F44").Value = 1201
$ws.Range("F45").Value = 1754
